$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (weather readings for Lisbon)
$newRows = @(
    @{ row = 14; values = @(12.44, 11.83, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:41:20 02-12-2025") },
    @{ row = 15; values = @(11.81, 11.21, 1019, 83, "few clouds", 20, "Lisbon", 18.504, 20, "19:52:54 02-12-2025") }
)

foreach ($entry in $newRows) {
    $r = $entry.row
    $vals = $entry.values
    for ($col = 1; $col -le $vals.Count; $col++) {
        $ws.Cells.Item($r, $col).Value = $vals[$col - 1]
    }
}
